$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.542273640632629
$ws.Range("B1").Value = 4.366642951965332
$ws.Range("C1").Value = 3.287420749664307
$ws.Range("D1").Value = 1.375802874565125
$ws.Range("E1").Value = 0.9475429058074951
